$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a single test row (row 3): reuse the existing "01-1234567-10" value in
# column C (same as row 2) and a new value "200000002C" in column D. The new
# D3 cell uses a wrapped-text style (matching the new cellXf added to the
# workbook for the source edit).
$ws.Range("C3").Value = "01-1234567-10"
$ws.Range("D3").Value = "200000002C"
$ws.Range("D3").WrapText = $true

# Leave the selection on the newly added cell, like in the authored edit.
$ws.Range("D3").Select() | Out-Null
